# Revert "update data export"
# This reverts the prior "update data export" commit by restoring the
# previous numeric values in columns C:L for the affected data rows of
# Sheet1 in sum_tagged_all_ops.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    "C2" = 11786696919.88891
    "D2" = 68707778478.51116
    "E2" = 80494475398.40007
    "F2" = 141778049092.0393
    "G2" = 123040385205.7217
    "H2" = 16587803787.86758
    "I2" = 2149860098.45
    "J2" = 10227981797.28891
    "K2" = 78935760275.80008
    "L2" = 1558715122.6
    "C3" = 2095495467.56
    "D3" = 11874474299.59
    "E3" = 13969969767.15
    "F3" = 17634714703.94
    "G3" = 14989507496.93
    "H3" = 1139045363.76
    "J3" = 788526165.84
    "K3" = 12663000465.43
    "C4" = 10323278549.5448
    "D4" = 8586738394.717028
    "E4" = 18910016944.26183
    "F4" = 29789212833.82301
    "G4" = 13391085667.5507
    "H4" = 4050922760.012311
    "I4" = 12347204406.26
    "J4" = 2376796728.554799
    "K4" = 10963535123.27183
    "L4" = 7946481820.99
    "C5" = 433704.8399999999
    "D5" = 2457660.7475
    "E5" = 2891365.5875
    "F5" = 2895764.4475
    "G5" = 2461399.78
    "H5" = 434364.6675
    "I5" = 0
    "J5" = 433704.84
    "K5" = 2891365.5875
    "C6" = 6981418187.888666
    "D6" = 39561369709.86514
    "E6" = 46542787897.7538
    "F6" = 72509301609.59196
    "G6" = 61632906368.13071
    "H6" = 10876395241.46125
    "J6" = 6981418187.888666
    "K6" = 46542787897.7538
    "L6" = -0.00000020302832126617428993647
    "C7" = 8449409756.46699
    "D7" = 46785239476.89042
    "E7" = 55234649233.35741
    "F7" = 83484737175.4998
    "G7" = 70056729236.9969
    "H7" = 10976071182.07291
    "I7" = 2451936756.43
    "J7" = 7831070276.74699
    "K7" = 54616309753.63741
    "L7" = 618339479.7200003
    "I8" = 712003308.25
    "C9" = 38302367911.9211
    "D9" = 31745404063.65495
    "E9" = 70047771975.57605
    "F9" = 125991790230.5634
    "G9" = 60331477127.82437
    "H9" = 1253846940.470587
    "I9" = 64406466162.26843
    "J9" = 356622884.0100001
    "K9" = 32102026947.66494
    "L9" = 37945745027.9111
    "C10" = 1299181050.89
    "E10" = 2356795765.86
    "F10" = 4430773073.11
    "I10" = 2425138210.78
    "L10" = 1299181050.89
    "C11" = 14491894433.37137
    "D11" = 10138046070.21396
    "E11" = 24629940503.58532
    "F11" = 56258741421.84528
    "G11" = 22671481222.51871
    "I11" = 33557291136.82656
    "K11" = 10138046070.21396
    "L11" = 14491894433.37137
    "C13" = 2233600213.363579
    "D13" = 2233570204.555042
    "E13" = 4467170417.918622
    "F13" = 7008425150.000095
    "G13" = 3504212574.197531
    "H13" = 2945860175.324294
    "I13" = 558352400.4782691
    "J13" = 1868029778.248265
    "K13" = 4101599982.803307
    "L13" = 365570435.1153149
    "G15" = 896703562.9
    "K15" = 873986284.6799999
    "L15" = 140594917.43
    "C16" = 786725549.6023514
    "D16" = 4458111417.312308
    "E16" = 5244836966.91466
    "F16" = 7185129211.911049
    "G16" = 6107359829.589524
    "H16" = 1077769382.321525
    "J16" = 786725549.6023514
    "K16" = 5244836966.91466
    "L16" = 0.000000026216730475135471140577
    "F17" = 0
    "G17" = 0
    "H17" = 0
    "C18" = 11723490376.01386
    "D18" = 48783070057.32639
    "E18" = 60506560433.34026
    "F18" = 88064258036.46448
    "G18" = 70797148996.39685
    "H18" = 15564019928.70356
    "I18" = 1703089111.364068
    "J18" = 10362905044.20415
    "K18" = 59145975101.53053
    "L18" = 1360585331.809717
    "C19" = 4439.488935185183
    "D19" = 17633.82587962963
    "E19" = 22073.31481481481
    "F19" = 28000617.61833333
    "G19" = 18622528.01861111
    "H19" = 9378089.599722221
    "J19" = 4439.488935185185
    "K19" = 22073.31481481481
    "C20" = 9052732403.779791
    "D20" = 42744683996.36911
    "E20" = 51797416400.1489
    "F20" = 73449196423.68405
    "G20" = 60423967776.9254
    "H20" = 12165030111.56466
    "I20" = 860198535.1939973
    "J20" = 8500745193.308492
    "K20" = 51245429189.6776
    "L20" = 551987210.4712989
    "F21" = 313089289.375
    "G21" = 252105556.19375
    "H21" = 60983733.18125
    "C22" = 5405662882.735476
    "D22" = 25345869977.45417
    "E22" = 30751532860.18965
    "F22" = 39312408397.66499
    "G22" = 31701535546.70824
    "H22" = 5435053992.325138
    "I22" = 2175818858.631612
    "J22" = 4431004403.013864
    "K22" = 29776874380.46804
    "L22" = 974658479.721612
    "C23" = 2551355021.62
    "D23" = 17502658082.76
    "E23" = 20054013104.38
    "F23" = 28793401667.05
    "G23" = 24584486245.37
    "H23" = 1022459827.08
    "I23" = 3186455594.6
    "J23" = 614924230.86
    "K23" = 18117582313.62
    "L23" = 1936430790.76
    "C24" = 7913285650.13
    "D24" = 6364840433.53
    "E24" = 14278126083.66
    "F24" = 22774859324.28
    "H24" = 12549681010.48
    "J24" = 7856238536.95
    "K24" = 14221078970.48
    "L24" = 57047113.17999994
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
